# Salt Lake Away Pass Types - add a "Match ID" column as the new column A.
# This shifts every existing column one place to the right (B..V -> C..W)
# and populates the new column A with the "Match ID" header/value (5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank column before column A; everything else
# (values, shared strings, styles, merged cells, column letters) shifts
# right by one automatically.
$ws.Columns("A").Insert() | Out-Null

# Header row (row 2 holds the column titles) gets the new "Match ID" label,
# bold like the other header cells (no border though - matches row style 3).
$ws.Range("A2").Value2 = "Match ID"
$ws.Range("A2:A19").Font.Bold = $true

# Rows 3 and 20 are hidden in this sheet. Writing straight into a hidden
# row causes this host to stamp an explicit row height on save, which the
# original edit never did - so briefly unhide, write, then re-hide.
$ws.Rows(3).Hidden = $false
$ws.Rows(20).Hidden = $false

# Data rows: every visible/hidden player row (4-20) gets Match ID = 5.
$ws.Range("A4:A20").Value2 = 5

$ws.Rows(20).Hidden = $true
$ws.Rows(3).Hidden = $true

# Restore the sheet's original selection, now anchored on the new column.
$ws.Range("A2:A19").Select() | Out-Null
